$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datasets")
$ws.Activate()
$ws.Range("F2").Value = "test"
